$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query (row 2, column B) previously returned an extra
# `Cohort` column sourced from an OPTIONAL MATCH on (co:cohort). That
# trailing RETURN-clause column is being dropped so the cases query only
# returns up through `Response to Treatment`.
$casesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['COTC007B','NCATS-COP01','GLIOMA01']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@
$casesQuery = $casesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $casesQuery

# Row heights were recomputed (by a newer Excel build) now that the cell
# content is a little shorter.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Selection moved from D4 to B2 and the view no longer pins a frozen/
# scrolled top-left cell.
[void]$ws.Range("B2").Select()
